# Auto update Excel log
# Appends a new sensor event row (row 54) to the "mmWave" log sheet,
# mirroring the most recent "PRESENCE_DETECTED" / "Active" entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mmWave")

$newRow = 54

# Column A holds a date-looking string ("2026-02-01"). Force the cell to
# Text format before assigning so Excel doesn't auto-convert it into a
# real date serial number, then restore the default "Normal" style so no
# stray number-format/style survives on the cell (matches the rest of the
# log, which stores these as plain text).
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "2026-02-01"
$ws.Cells.Item($newRow, 1).Style = "Normal"

$ws.Cells.Item($newRow, 2).Value = "16:00:49"
$ws.Cells.Item($newRow, 3).Value = "16:00"
$ws.Cells.Item($newRow, 4).Value = "Living Room"
$ws.Cells.Item($newRow, 5).Value = "PRESENCE_DETECTED"
$ws.Cells.Item($newRow, 6).Value = "Active"
